# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# OFF sheet - row 3 (R / Road) updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 416
$wsOff.Range("C3").Value = 294
$wsOff.Range("D3").Value = 93
$wsOff.Range("E3").Value = 51

# DEF sheet - row 3 (R / Road) updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 486
$wsDef.Range("C3").Value = 361
$wsDef.Range("D3").Value = 122
$wsDef.Range("E3").Value = 54
$wsDef.Range("F3").Value = 4
$wsDef.Range("G3").Value = 6
